$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Window was moved/resized (no longer off the top of the screen) between edits
$win = $excel.ActiveWindow
$win.Top = 0
$win.Height = 14200

# Rename the sheet (Sheet1 -> Q4)
$ws.Name = "Q4"

# New header row (D1:F1)
$ws.Range("D1").Value = "total memory accesses"
$ws.Range("E1").Value = "hit"
$ws.Range("F1").Value = "miss"

# New column D width, matching the style already used for column B
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Row 2 - total memory accesses / hit / miss
$ws.Range("D2").Value = 49642128
$ws.Range("E2").Value = 35991042
$ws.Range("F2").Value = 13651086

# Row 3
$ws.Range("E3").Value = 38469661
$ws.Range("F3").Value = 11172467

# Row 4
$ws.Range("E4").Value = 40482042
$ws.Range("F4").Value = 9160086

# Row 5
$ws.Range("E5").Value = 45397029
$ws.Range("F5").Value = 4245099

# Row 6
$ws.Range("E6").Value = 46488794
$ws.Range("F6").Value = 3153334

# Row 7
$ws.Range("E7").Value = 47406881
$ws.Range("F7").Value = 2235247

# Row 8
$ws.Range("E8").Value = 48835377
$ws.Range("F8").Value = 806751

# Row 9
$ws.Range("E9").Value = 49046540
$ws.Range("F9").Value = 595588

# Row 10
$ws.Range("E10").Value = 49231283
$ws.Range("F10").Value = 410845

# Update the active selection to match the post-edit cursor position
[void]$ws.Range("F16").Select()
